$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replaces the previously-empty row 4, shifted down to row 5
$ws.Range("A5").Value = "Dr Riaz Khan"
$ws.Range("B5").Value = "riaz.khan.ruet@gmail.com"
$ws.Range("C5").Value = "Machine Learning, Deep Learning"
$ws.Range("D5").Value = "sites\google.com\riaz-mte-16"

# Row 6: new row with similar data
$ws.Range("A6").Value = "Dr Riaz Khan"
$ws.Range("B6").Value = "riaz.khan.ruet@gmail.com"
$ws.Range("C6").Value = "Machine Learning, Deep Learning"
$ws.Range("D6").Value = "sites\google.com\riaz-mte-16 f"

# Clear out old row 4 (which previously held empty inline-string cells)
$ws.Range("A4:D4").ClearContents()
